$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$f2Text = @'
日语/Provide a list of requested number of items.  任务目标: 设计不同人物对相关主题的看法，先用中高级日语回答，然后提供中文翻译。然后对重点单词或语法做说明。  Number of items: 10  Example output item: ''' item x: [日语句子][中文翻译][重点单词或语法说明] ''' Reply in the following format:     - item 1     - item 2     - item 3
'@

$f12Text = @'
专家发言/提供20个各有特长的专家做出的专业解答，思维方式如下:###对每个子问题构建一系列清晰的思维步骤，记录并解释每个推理过程。
确保思维链条完整、严谨，避免遗漏任何关键逻辑环节。
假设验证：
对每个可能的解决方案进行假设测试，分析其合理性和潜在局限性。
根据反馈修正错误，确保假设的准确性和可行性。
多策略思考：
如果现有方法未能得出理想结果，迅速切换策略，尝试不同路径，探索备选方案。
综合多个角度分析问题，形成更全面的解决方案。
推理总结：
汇总所有推理步骤，结合多种思考路径，最终得出一个具有说服力的答案。
阐明每个结论背后的逻辑及推理依据，确保解释条理清晰、严密。###
回答样式:[解答方法][具体解答]...Final output are in the following format:     - 段落 1     - 段落 2     - 段落 3

'@

$ws.Range("F2").Value = $f2Text
$ws.Range("F12").Value = $f12Text
$ws.Range("F13").ClearContents()

$ws.Rows.Item(12).RowHeight = 374.4
$ws.Rows.Item(13).EntireRow.AutoFit()

$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F13").Select()
